$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.031.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.66%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.884.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.12%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.7376"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.48%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'241.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.06%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3160"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.17%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07172"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.96%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -2.39%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.08317"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.19%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.7555"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.40%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.928.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.09%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.403"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.78%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'92.51"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'6.146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.30%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'30.057.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.84%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'249.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.83%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -1.05%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.000007850"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.36%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'2.157.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.88%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.05%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'7.888"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.16%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.1569"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.84%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'9.259"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.05%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'164.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.56%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'18.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.23%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.043"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.72%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.472"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.36%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'4.546"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.00%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.531"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.10%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'4.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.56%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.05318"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.83%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.72%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.7662"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.93%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.9992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.29%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'2.728"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.70%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01956"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.80%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.760"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.39%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.4552"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.06%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.8792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.32%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'6.036"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.07%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'Maker"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.086.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.32%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'72.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.29%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'104.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.77%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.07%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -0.33%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'7.528"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.33%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'9.536"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.99%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'2.061.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.00%  "
$ws.Range("E51").Style = "Normal"
